$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E (H.H. horas hombre) : header fill + per-row time values ---
# Header E1 gets a yellow fill (creates style index 5: General + yellow fill)
$ws.Range("E1").Interior.Color = 65535

# Data rows E2:E13 -> time-of-day values (duration in minutes / 1440) with
# time number format "h:mm" and the same yellow fill (creates style index 6:
# time format + yellow fill). Values match column F where the same duration
# was already tracked, except E8/E9/E11/E12/E13 which differ slightly.
$ws.Range("E2").Value = 15/1440
$ws.Range("E2").NumberFormat = "h:mm"
$ws.Range("E2").Interior.Color = 65535

$ws.Range("E3").Value = 30/1440
$ws.Range("E3").NumberFormat = "h:mm"
$ws.Range("E3").Interior.Color = 65535

$ws.Range("E4").Value = 15/1440
$ws.Range("E4").NumberFormat = "h:mm"
$ws.Range("E4").Interior.Color = 65535

$ws.Range("E5").Value = 30/1440
$ws.Range("E5").NumberFormat = "h:mm"
$ws.Range("E5").Interior.Color = 65535

$ws.Range("E6").Value = 30/1440
$ws.Range("E6").NumberFormat = "h:mm"
$ws.Range("E6").Interior.Color = 65535

$ws.Range("E7").Value = 15/1440
$ws.Range("E7").NumberFormat = "h:mm"
$ws.Range("E7").Interior.Color = 65535

$ws.Range("E8").Value = 20/1440
$ws.Range("E8").NumberFormat = "h:mm"
$ws.Range("E8").Interior.Color = 65535

$ws.Range("E9").Value = 75/1440
$ws.Range("E9").NumberFormat = "h:mm"
$ws.Range("E9").Interior.Color = 65535

$ws.Range("E10").Value = 30/1440
$ws.Range("E10").NumberFormat = "h:mm"
$ws.Range("E10").Interior.Color = 65535

$ws.Range("E11").Value = 20/1440
$ws.Range("E11").NumberFormat = "h:mm"
$ws.Range("E11").Interior.Color = 65535

# --- New rows 12 and 13 ---
$ws.Range("C12").Value = "Limpieza repuestos Nuevos"
$ws.Range("C12").Interior.Color = 65535
$ws.Range("E12").Value = 35/1440
$ws.Range("E12").NumberFormat = "h:mm"
$ws.Range("E12").Interior.Color = 65535

$ws.Range("C13").Value = "Embalar Repuestos en desuso"
$ws.Range("C13").Interior.Color = 65535
$ws.Range("E13").Value = 30/1440
$ws.Range("E13").NumberFormat = "h:mm"
$ws.Range("E13").Interior.Color = 65535

# --- New column G (Cargos) ---
$ws.Range("G1").Value = "Cargos"

# Create style 7 (centered, no fill, General) on a blank cell first
$ws.Range("G3").HorizontalAlignment = -4108

# G2 / G9 carry an actual "cargo" duration (centered + time format -> style 8)
$ws.Range("G2").Value = 270/1440
$ws.Range("G2").NumberFormat = "h:mm"
$ws.Range("G2").HorizontalAlignment = -4108

$ws.Range("G9").Value = 270/1440
$ws.Range("G9").NumberFormat = "h:mm"
$ws.Range("G9").HorizontalAlignment = -4108

# Remaining G cells in both groups stay blank, just centered (style 7)
$ws.Range("G4").HorizontalAlignment = -4108
$ws.Range("G5").HorizontalAlignment = -4108
$ws.Range("G6").HorizontalAlignment = -4108
$ws.Range("G7").HorizontalAlignment = -4108
$ws.Range("G8").HorizontalAlignment = -4108
$ws.Range("G10").HorizontalAlignment = -4108
$ws.Range("G11").HorizontalAlignment = -4108
$ws.Range("G12").HorizontalAlignment = -4108
$ws.Range("G13").HorizontalAlignment = -4108

# Merge the two "Cargos" groups
$ws.Range("G2:G8").Merge()
$ws.Range("G9:G13").Merge()

# --- Totals row 14 ---
$ws.Range("E14").Formula = "=SUM(E2:E13)"
$ws.Range("E14").NumberFormat = "h:mm"
$ws.Range("E14").Interior.Color = 65535

$ws.Range("F14").Formula = "=SUM(F2:F13)"
$ws.Range("F14").NumberFormat = "h:mm"
$ws.Range("F14").Interior.Color = 65535

$ws.Range("G14").Formula = "=SUM(G2:G13)"
$ws.Range("G14").NumberFormat = "h:mm"
$ws.Range("G14").Interior.Color = 65535

# --- Row heights / column widths ---
$ws.Rows.Item(11).RowHeight = 45
$ws.Columns.Item(3).ColumnWidth = 42.5
$ws.Columns.Item(5).ColumnWidth = 9.140625

# --- Final selection state ---
$ws.Range("E14:G14").Select()
$excel.ActiveCell = $ws.Range("G14")
